# Update the Campaign sheet: "HL Sub Group" value (E2) changes from "BAS" to "None"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Campaign")
$ws.Range("E2").Value = "None"
